$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2:D5 values (Power column)
$ws.Range("D2").Value = 0.136
$ws.Range("D3").Value = 0.185
$ws.Range("D4").Value = 0.195
$ws.Range("D5").Value = 0.07

# D2:D3 already used a custom numeric display; apply the same "0.000" format
# consistently across D2:D5 so they all share the one format.
$ws.Range("D2:D5").NumberFormat = "0.000"

# Add the Total row formula for D6, styled like the adjacent B6/C6 totals
# (bold, centered, General number format) by copying C6's format over.
$ws.Range("D6").Formula = "=D2+D3+D4+D5"
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)

# Move the active selection to D7
$ws.Range("D7").Select()
